# Generate Report for Handoff
# The translation run moved from "In Translation" to "Ready for handoff"; the
# Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps were
# refreshed, and the Status/zh-cn/de-de columns widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" --------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps -------------------------------------------------
# Overview!G2 (Latest HO Xliff Generate Date) and de-de!H2 (Latest Handoff
# Datetime) both held "2016-09-03 04:42:45" and now hold "...04:43:36".
$overview.Range("G2").Value = "2016-09-03 04:43:36"
$dede.Range("H2").Value     = "2016-09-03 04:43:36"

# zh-cn!H2 (Latest Handoff Datetime) held "2016-09-03 04:42:40" and now
# holds "...04:43:31".
$zhcn.Range("H2").Value = "2016-09-03 04:43:31"

# --- Widen the columns that now hold the longer "Ready for handoff" text -
$overview.Range("E2").ColumnWidth = 16.25
$overview.Range("F2").ColumnWidth = 16.25
$zhcn.Range("C2").ColumnWidth     = 16.25
$dede.Range("C2").ColumnWidth     = 16.25
